$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 49068
$ws.Range("C2").Value = 24
$ws.Range("E2").Value = 0.02105258878195004
$ws.Range("F2").Value = 42.0303211156045
$ws.Range("G2").Value = 0.04241296643536585
$ws.Range("I2").Value = 44752
$ws.Range("J2").Value = 4

# Row 3
$ws.Range("A3").Value = 42313
$ws.Range("B3").Value = 34
$ws.Range("C3").Value = 59
$ws.Range("D3").Value = 0.02342187848
$ws.Range("E3").Value = 0.02419069819893193
$ws.Range("F3").Value = 32.53106876716668
$ws.Range("G3").Value = 0.03282485303595206
$ws.Range("H3").Value = 44762
$ws.Range("I3").Value = 44787
$ws.Range("J3").Value = 25

# Row 4
$ws.Range("A4").Value = 49074
$ws.Range("E4").Value = 0.0208499967650379
$ws.Range("F4").Value = 32.20251313958988
$ws.Range("G4").Value = 0.03249528827451198

# Row 5
$ws.Range("A5").Value = 57932
$ws.Range("B5").Value = 174
$ws.Range("C5").Value = 215
$ws.Range("D5").Value = 0.01710604896
$ws.Range("E5").Value = 0.02114152160714397
$ws.Range("F5").Value = 233.7830013943446
$ws.Range("G5").Value = 0.2359091018960799
$ws.Range("H5").Value = 44902
$ws.Range("I5").Value = 44943
$ws.Range("J5").Value = 41

# Row 6
$ws.Range("A6").Value = 46842
$ws.Range("B6").Value = 216
$ws.Range("C6").Value = 222
$ws.Range("D6").Value = 0.02115342229
$ws.Range("E6").Value = 0.02287997481358652
$ws.Range("F6").Value = 80.87517330983997
$ws.Range("G6").Value = 0.08162048201546712
$ws.Range("H6").Value = 44944
$ws.Range("I6").Value = 44950
$ws.Range("J6").Value = 6

# Row 7 (previously the last row, values change + C7 newly present)
$ws.Range("A7").Value = 43741
$ws.Range("B7").Value = 223
$ws.Range("C7").Value = 232
$ws.Range("D7").Value = 0.02265457194
$ws.Range("E7").Value = 0.02324124188836351
$ws.Range("F7").Value = 25.66153021136832
$ws.Range("G7").Value = 0.0258963157598957
$ws.Range("H7").Value = 44951
$ws.Range("I7").Value = 44960
$ws.Range("J7").Value = 9

# New Row 8
$ws.Range("A8").Value = 44594
$ws.Range("B8").Value = 244
$ws.Range("C8").Value = 249
$ws.Range("D8").Value = 0.02222203984
$ws.Range("E8").Value = 0.02408751406772674
$ws.Range("F8").Value = 83.18895771124652
$ws.Range("G8").Value = 0.08394702921776176
$ws.Range("H8").Value = 44972
$ws.Range("H8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I8").Value = 44977
$ws.Range("I8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J8").Value = 5
$ws.Range("J8").NumberFormat = "0"

# New Row 9
$ws.Range("A9").Value = 40487
$ws.Range("B9").Value = 251
$ws.Range("C9").Value = 267
$ws.Range("D9").Value = 0.02447512067
$ws.Range("E9").Value = 0.02036221
$ws.Range("F9").Value = -166.5194142962899
$ws.Range("G9").Value = -0.1680445512589989
$ws.Range("H9").Value = 44979
$ws.Range("H9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I9").Value = 44995
$ws.Range("I9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J9").Value = 16
$ws.Range("J9").NumberFormat = "0"
